$wb = $excel.ActiveWorkbook

$wsConfigs = $wb.Worksheets.Item("configs")
$wsStudents = $wb.Worksheets.Item("students")

# Leave a "last selection" of B7 behind on the configs sheet (matches author's trace)
$wsConfigs.Activate()
$wsConfigs.Range("B7").Select()

# Switch to the students sheet and append the two new students to the bottom of the table
$wsStudents.Activate()
$wsStudents.Range("D6").Select()

$wsStudents.Range("A5").Value = "Younes"
$wsStudents.Range("B5").Value = "Sayeh"
$wsStudents.Range("A6").Value = "Philipona"
$wsStudents.Range("B6").Value = "Sylvain"

$wsStudents.Range("D7").Select()
